$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 8694
$ws.Range("C3:C13").Value = 7769
$ws.Range("C14:C28").Value = 7312
$ws.Range("C68:C252").Value = 7310
